$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) The top data row (row 7) gains a "Total" quantity value in P7.
#    Previously P7 was blank; it now shows the text "2.0000". The
#    cell's number format is numeric (0.00), so Excel would normally
#    coerce the text into a real number - flip to Text, assign, then
#    restore the original display format to keep the value textual.
# ------------------------------------------------------------------
$p7fmt = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "2.0000"
$ws.Range("P7").NumberFormat = $p7fmt

# ------------------------------------------------------------------
# 2) Insert a brand-new data row above the old totals row (row 8),
#    for a second item ("سرنجات 5 سم"). This pushes the previous
#    row 8 (totals) down to row 9, and the previous row 9 (footer)
#    down to row 10 - exactly mirroring row 7's own layout/merges.
# ------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

# Clone formatting from row 7 into the freshly inserted row 8 so the
# new item row looks identical in style to the first item row.
$ws.Range("A7:O7").Copy()
$ws.Range("A8:O8").PasteSpecial(-4122)
$ws.Range("P7:Q7").Copy()
$ws.Range("P8:Q8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(8).RowHeight = 24.75

# Re-create the merges for the new row 8 (mirrors row 7's merges).
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# Populate the new row 8 with the second item's data. A8 is a real
# number (item index); the rest are text labels. L8 and P8 sit on
# numeric-looking display formats, so they need the same Text-format
# round-trip trick as P7 above to avoid being coerced into numbers.
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "سرنجات 5 سم"
$ws.Range("H8").Value = "0:0"

$l8fmt = $ws.Range("L8").NumberFormat
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "0"
$ws.Range("L8").NumberFormat = $l8fmt

$ws.Range("N8").Value = "3.00"

$p8fmt = $ws.Range("P8").NumberFormat
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "3.0000"
$ws.Range("P8").NumberFormat = $p8fmt

$ws.Range("Q8").Value = "1:0"

# ------------------------------------------------------------------
# 3) The totals row (old row 8, now row 9) changes shape: it used to
#    be a P:Q merge holding 0; now it is an N:Q merge holding 5, with
#    a taller row and a smaller font.
# ------------------------------------------------------------------
$ws.Range("P9:Q9").UnMerge()
$ws.Rows.Item(9).RowHeight = 26.25
$ws.Range("N9:Q9").Merge()
$ws.Range("N9").Value = 5
$ws.Range("N9:Q9").Font.Size = 13

# ------------------------------------------------------------------
# 4) The footer row (old row 9, now row 10) simply gets an updated
#    timestamp reflecting the new upload time.
# ------------------------------------------------------------------
$ws.Range("A10").Value = "Wednesday, 17 September, 2025 10:44 PM"

Write-Output "edit applied"
